$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the StatQuery column (C2:C4) which holds the filter Cypher query.
# The "breed" IN ['V'] filter is removed (reset to empty) and the
# "stage_of_disease" IN [] filter is set to ['V'] instead.
$oldQuery = $ws.Range("C2").Value2
$newQuery = $oldQuery.Replace("AND (size([`'V`']) = 0 OR demo.breed IN [`'V`'])", "AND (size([]) = 0 OR demo.breed IN [])")
$newQuery = $newQuery.Replace("AND (size([]) = 0 OR diag.stage_of_disease IN [])", "AND (size([`'V`']) = 0 OR diag.stage_of_disease IN [`'V`'])")

$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# Zoom out the sheet view from 55% to 40%.
$ws.Application.ActiveWindow.Zoom = 40

# Adjust the workbook window position/size (best effort; some hosts may not
# persist these window-chrome settings into the saved file).
$win = $wb.Windows.Item(1)
$win.Left = 28680
$win.Top = -105
$win.Width = 29040
$win.Height = 15840
